# COCO3_Comp.xlsx edit
#
# Commit message summary: cart list updated (FE00 addressing fixes for
# PITFALL / RAMPAGE, bank-switching rework). The workbook-level data change
# that shows up in the sheet is a new cart entry, "Popstar Demo", added to
# the compatibility table as row 31 (pushing the three "known issue" notes
# below it down from rows 33-35 to rows 34-36), plus the corresponding
# sheet-view selection update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 31. Excel shifts every row at/after 31
# down by one (the former rows 33/34/35 become 34/35/36) and extends the
# sheet's used range/dimension automatically.
$ws.Rows.Item(31).Insert()

# Fill in the new cart entry: name in column A, "Y" (pass) in the D and E
# load-test columns, matching the style already used by the surrounding
# rows (centered text, picked up automatically from the row below on
# insert).
$ws.Range("A31").Value = "Popstar Demo"
$ws.Range("D31").Value = "Y"
$ws.Range("E31").Value = "Y"

# Reflect the author's on-screen state when they saved: scrolled down so
# row 4 is at the top of the window, with the new row's first cell
# selected.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A31").Select()
